$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 37; $row++) {
    $hCell = $ws.Cells.Item($row, 8)
    $hVal = $hCell.Value()
    $hCell.Value = $hVal - 1

    $iCell = $ws.Cells.Item($row, 9)
    $origFormat = $iCell.NumberFormat()
    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"
    $iCell.NumberFormat = $origFormat
}
